$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell B1: "Description" -> "ICE Dimension"
$ws.Range("B1").Value = "ICE Dimension"

# Move selection / active cell to C8 (was a full-column selection A1:A8)
$ws.Range("C8").Select()
